# Auto-generated edit script applying Shinryu_Profits market data updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 135.96428
$ws.Range("I53").Value = 126.4375
$ws.Range("J53").Value = 148.66667
$ws.Range("K53").Value = 126.4375
$ws.Range("L53").Value = 148.66667
$ws.Range("M53").Value = 510.5625
$ws.Range("N53").Value = -1422.66667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 38845.7
$ws.Range("J7").Value = 38845.7
$ws.Range("L7").Value = 38845.7
$ws.Range("N7").Value = -39073.7

$ws.Range("H52").Value = 33700
$ws.Range("J52").Value = 33700
$ws.Range("L52").Value = 33700
$ws.Range("N52").Value = -34336

$ws.Range("H111").Value = 28600
$ws.Range("J111").Value = 28600
$ws.Range("L111").Value = 28600
$ws.Range("N111").Value = -36780

$ws.Range("H127").Value = 34325
$ws.Range("J127").Value = 34325
$ws.Range("L127").Value = 34325
$ws.Range("N127").Value = -44245

$ws.Range("H129").Value = 35582.145
$ws.Range("J129").Value = 35582.145
$ws.Range("L129").Value = 35582.145
$ws.Range("N129").Value = -45582.145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 22000
$ws.Range("J6").Value = 22000
$ws.Range("L6").Value = 22000
$ws.Range("N6").Value = -22226

$ws.Range("H51").Value = 34933.332
$ws.Range("J51").Value = 34933.332
$ws.Range("L51").Value = 34933.332
$ws.Range("N51").Value = -35915.332

$ws.Range("H52").Value = 32933.332
$ws.Range("J52").Value = 32933.332
$ws.Range("L52").Value = 32933.332
$ws.Range("N52").Value = -33459.332

$ws.Range("H55").Value = 38000
$ws.Range("J55").Value = 38000
$ws.Range("L55").Value = 38000
$ws.Range("N55").Value = -38546

$ws.Range("H86").Value = 2966.6667
$ws.Range("I86").Value = 1950
$ws.Range("K86").Value = 1950
$ws.Range("M86").Value = -827

$ws.Range("H89").Value = 2966.6667
$ws.Range("I89").Value = 1950
$ws.Range("K89").Value = 9750
$ws.Range("M89").Value = -4134

$ws.Range("H105").Value = 2591.31
$ws.Range("I105").Value = 1291.1177
$ws.Range("J105").Value = 2857.6145
$ws.Range("K105").Value = 1291.1177
$ws.Range("L105").Value = 2857.6145
$ws.Range("M105").Value = 455.8823
$ws.Range("N105").Value = -6351.6145

$ws.Range("H107").Value = 2192.2
$ws.Range("I107").Value = 2192.2
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2192.2
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -272.1999999999998
$ws.Range("N107").Value = $null

$ws.Range("H118").Value = 36000
$ws.Range("J118").Value = 36000
$ws.Range("L118").Value = 36000
$ws.Range("N118").Value = -39314

$ws.Range("H121").Value = 32933.332
$ws.Range("J121").Value = 32933.332
$ws.Range("L121").Value = 32933.332
$ws.Range("N121").Value = -36427.332

$ws.Range("H127").Value = 19337.5
$ws.Range("I127").Value = 17950
$ws.Range("K127").Value = 17950
$ws.Range("M127").Value = -12990

$ws.Range("H134").Value = 1666.5385
$ws.Range("I134").Value = 1666.5385
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4999.6155
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2464.6155
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 30000
$ws.Range("J18").Value = 30000
$ws.Range("L18").Value = 30000
$ws.Range("N18").Value = -30460

$ws.Range("H31").Value = 1837.3334
$ws.Range("I31").Value = 2038.5
$ws.Range("J31").Value = 1722.381
$ws.Range("K31").Value = 2038.5
$ws.Range("L31").Value = 1722.381
$ws.Range("M31").Value = -1743.5
$ws.Range("N31").Value = -2312.381

$ws.Range("H34").Value = 1837.3334
$ws.Range("I34").Value = 2038.5
$ws.Range("J34").Value = 1722.381
$ws.Range("K34").Value = 2038.5
$ws.Range("L34").Value = 1722.381
$ws.Range("M34").Value = -1836.5
$ws.Range("N34").Value = -2126.381

$ws.Range("H99").Value = 2099.7334
$ws.Range("I99").Value = 1629.0667
$ws.Range("J99").Value = 2570.4
$ws.Range("K99").Value = 1629.0667
$ws.Range("L99").Value = 2570.4
$ws.Range("M99").Value = -131.0667000000001
$ws.Range("N99").Value = -5566.4

$ws.Range("H105").Value = 478116.34
$ws.Range("I105").Value = 715895
$ws.Range("J105").Value = 2559
$ws.Range("K105").Value = 715895
$ws.Range("L105").Value = 2559
$ws.Range("M105").Value = -714148
$ws.Range("N105").Value = -6053

$ws.Range("H107").Value = 500891.5
$ws.Range("I107").Value = 1000790.1
$ws.Range("J107").Value = 992.9
$ws.Range("K107").Value = 1000790.1
$ws.Range("L107").Value = 992.9
$ws.Range("M107").Value = -998870.1
$ws.Range("N107").Value = -4832.9

$ws.Range("H114").Value = 31950
$ws.Range("J114").Value = 31950
$ws.Range("L114").Value = 31950
$ws.Range("N114").Value = -40628

$ws.Range("H117").Value = 46106.855
$ws.Range("J117").Value = 46106.855
$ws.Range("L117").Value = 46106.855
$ws.Range("N117").Value = -55284.855

$ws.Range("H123").Value = 35000
$ws.Range("J123").Value = 35000
$ws.Range("L123").Value = 35000
$ws.Range("N123").Value = -44800

$ws.Range("H126").Value = 2099.7334
$ws.Range("I126").Value = 1629.0667
$ws.Range("J126").Value = 2570.4
$ws.Range("K126").Value = 4887.2001
$ws.Range("L126").Value = 7711.200000000001
$ws.Range("M126").Value = -2417.2001
$ws.Range("N126").Value = -12651.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 32.95238
$ws.Range("I12").Value = 19.25
$ws.Range("J12").Value = 36.17647
$ws.Range("K12").Value = 57.75
$ws.Range("L12").Value = 108.52941
$ws.Range("M12").Value = 115.25
$ws.Range("N12").Value = -454.52941

$ws.Range("H40").Value = 142
$ws.Range("I40").Value = 102.5
$ws.Range("J40").Value = 300
$ws.Range("K40").Value = 410
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -341
$ws.Range("N40").Value = -1338

$ws.Range("H75").Value = 5878.25
$ws.Range("I75").Value = 675.3333
$ws.Range("J75").Value = 9000
$ws.Range("K75").Value = 2025.9999
$ws.Range("L75").Value = 27000
$ws.Range("M75").Value = -1027.9999
$ws.Range("N75").Value = -28996

$ws.Range("H78").Value = 5878.25
$ws.Range("I78").Value = 675.3333
$ws.Range("J78").Value = 9000
$ws.Range("K78").Value = 6077.9997
$ws.Range("L78").Value = 81000
$ws.Range("M78").Value = -1085.9997
$ws.Range("N78").Value = -90984

$ws.Range("H107").Value = 24391038
$ws.Range("I107").Value = 199.57143
$ws.Range("J107").Value = 37038140
$ws.Range("K107").Value = 598.71429
$ws.Range("L107").Value = 111114420
$ws.Range("M107").Value = 1321.28571
$ws.Range("N107").Value = -111118260

$ws.Range("H117").Value = 2393.5293
$ws.Range("I117").Value = 1591.7778
$ws.Range("J117").Value = 3295.5
$ws.Range("K117").Value = 4775.3334
$ws.Range("L117").Value = 9886.5
$ws.Range("M117").Value = -1333.3334
$ws.Range("N117").Value = -16770.5

$ws.Range("H131").Value = 2539.4126
$ws.Range("J131").Value = 2811.8909
$ws.Range("L131").Value = 8435.672699999999
$ws.Range("N131").Value = -18515.6727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 30000
$ws.Range("J32").Value = 30000
$ws.Range("L32").Value = 30000
$ws.Range("N32").Value = -30592

$ws.Range("H98").Value = 12643
$ws.Range("J98").Value = 12643
$ws.Range("L98").Value = 12643
$ws.Range("N98").Value = -18633

$ws.Range("H132").Value = 13666.333
$ws.Range("I132").Value = 15857
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 47571
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -45041
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 25500
$ws.Range("J108").Value = 25500
$ws.Range("L108").Value = 25500
$ws.Range("N108").Value = -33180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 20501
$ws.Range("I3").Value = 7000
$ws.Range("J3").Value = 25001.334
$ws.Range("K3").Value = 7000
$ws.Range("L3").Value = 25001.334
$ws.Range("M3").Value = -6886
$ws.Range("N3").Value = -25229.334

$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = $null

$ws.Range("H9").Value = 70007
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 70007
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 70007
$ws.Range("M9").Value = $null
$ws.Range("N9").Value = -70287

$ws.Range("H10").Value = 70006
$ws.Range("J10").Value = 70006
$ws.Range("L10").Value = 70006
$ws.Range("N10").Value = -70344

$ws.Range("H12").Value = 250052510
$ws.Range("I12").Value = 1000000000
$ws.Range("K12").Value = 1000000000
$ws.Range("M12").Value = -999999858

$ws.Range("H13").Value = 70000
$ws.Range("I13").Value = 70000
$ws.Range("K13").Value = 70000
$ws.Range("M13").Value = -69860

$ws.Range("H107").Value = 675.3333
$ws.Range("I107").Value = 675.3333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2025.9999
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -105.9999
$ws.Range("N107").Value = $null
